$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"

$ws.Cells.Item($row, 4).Value = 44911
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100103
$ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value = 100103003
$ws.Cells.Item($row, 10).Value = "Damasco"
$ws.Cells.Item($row, 11).Value = "Castle Brite"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 600
$ws.Cells.Item($row, 14).Value = 19000
$ws.Cells.Item($row, 15).Value = 20000
$ws.Cells.Item($row, 16).Value = 19500
$ws.Cells.Item($row, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 1219
$ws.Cells.Item($row, 20).Value = 16
